$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) rows 2-11 from 45221 to 45224
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
